# Swap the order of the two comma-separated "Recorded By" values in
# column G for the specific known combinations that changed order.
# This mirrors an upstream sync where "X, Y" became "Y, X" for three
# exact value pairs:
#   "dnasr281@gmail.com, System"          -> "System, dnasr281@gmail.com"
#   "System, backup@backdoor.com"         -> "backup@backdoor.com, System"
#   "dnasr281@gmail.com, admin@admin.com" -> "admin@admin.com, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

$colIndex = 7  # Column G ("Recorded By")

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $colIndex)
    $val = $cell.Value2

    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value2 = "System, dnasr281@gmail.com"
    }
    elseif ($val -eq "System, backup@backdoor.com") {
        $cell.Value2 = "backup@backdoor.com, System"
    }
    elseif ($val -eq "dnasr281@gmail.com, admin@admin.com") {
        $cell.Value2 = "admin@admin.com, dnasr281@gmail.com"
    }
}

$wb.Save()
